$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01838528797222151
$ws.Range("C2").Value = 0.1943756155675498

$ws.Range("B3").Value = 0.07300731111934443
$ws.Range("C3").Value = 0.240242512553099

$ws.Range("B4").Value = 0.8938796618671565
$ws.Range("C4").Value = 0.5191793218823673

$ws.Range("B5").Value = 0.9942781507963687
$ws.Range("C5").Value = 0.4637003269448957

$ws.Range("B6").Value = 0.9914961385628619
$ws.Range("C6").Value = 0.7904240810945034

$ws.Range("B7").Value = 0.9674376379532397
$ws.Range("C7").Value = 0.3622432403944318

$ws.Range("B8").Value = 0.006631307601928711
$ws.Range("C8").Value = 0.1713283348083496
